# Update "paises.xlsx" (Pais sheet) with the latest COVID-19 country stats.
# Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $name, $b, $c, $d, $e, $f, $g, $h) {
    if ($name -ne $null) {
        $ws.Cells.Item($r, 1).Value = $name
    }
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 14:20"

# --- Rows whose country stays in place, only the counters change ---

Set-Row 4 $null 6637899 1652 3918494 2521951 0 33 197454      # Estados Unidos
Set-Row 5 $null 4665469 8090 3626106 961815 0 42 77548        # India
Set-Row 38 $null 94211 736 84404 9249 0 1 558                 # Kuwait
Set-Row 58 $null 53120 1202 37524 15260 0 14 336               # Nepal
Set-Row 82 $null 19557 341 16247 2680 0 1 630                  # Dinamarca
Set-Row 86 $null 15694 140 13128 1920 0 4 646                  # Republica de Macedonia
Set-Row 87 $null 14237 44 10373 3569 0 2 295                   # Senegal
Set-Row 134 $null 3172 3 2983 177 0 0 12                       # Sri Lanka

# --- Rows whose ranking swaps because the updated totals cross each other ---

# Japon / Bielorrusia (rows 47-48): Bielorrusia's new total overtakes Japon's.
Set-Row 47 "Bielorrusia" 73975 191 72547 684 0 6 744
Set-Row 48 "Japon" 73901 0 65590 6899 0 0 1412

# Congo / Hong Kong (rows 115-116): Hong Kong's new total overtakes Congo's.
Set-Row 115 "Hong Kong" 4939 13 4613 226 0 1 100
Set-Row 116 "Congo" 4928 0 3887 953 0 0 88

# Somalia / Mayotte / Gambia (rows 129-131): Gambia's new total ties/overtakes
# Somalia's, and both now rank above Mayotte.
Set-Row 129 "Gambia" 3376 14 1617 1657 0 2 102
Set-Row 130 "Somalia" 3376 0 2791 487 0 0 98
Set-Row 131 "Mayotte" 3374 0 2964 370 0 0 40
